$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The workbook originally has 4 sheets:
#   "Sheet4" - contains the actual Team Offense data table (A1:Y143)
#   "Sheet1" - empty placeholder sheet
#   "Sheet2" - empty placeholder sheet
#   "Sheet3" - empty placeholder sheet
#
# Restructure the workbook data model: remove the two unused empty
# placeholder sheets ("Sheet2" and "Sheet3"), then rename the data sheet
# ("Sheet4") to "Sheet2" so the remaining sheets are "Sheet2" (data) and
# "Sheet1" (empty).

$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$wsData = $wb.Worksheets.Item("Sheet4")
$wsData.Activate()
$wsData.Range("A1").Select()
$wsData.Name = "Sheet2"
